$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-22 23:17:48"
$wsZhCn.Range("H2").Value = "2016-03-22 23:18:11"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-22 23:17:52"
$wsDeDe.Range("H2").Value = "2016-03-22 23:18:17"
